# Append a new row (Wordpress job) at row 4, shifting the two existing
# rows (FX arbitrage, AWS/HP) down by one, and refresh the "fetched at"
# timestamp (column A) for every data row to 2025-12-28 12:36:26.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2025-12-28 12:36:26"

# 1) Refresh timestamps on the two rows that are NOT moving (rows 2 & 3).
$ws.Range("A2").Value = $newTimestamp
$ws.Range("A3").Value = $newTimestamp

# 2) Insert a new row above the current row 4 - this pushes the existing
#    row 4 (FX latency arbitrage) down to row 5, and row 5 (AWS/HP) down
#    to row 6.
$ws.Rows.Item(4).Insert()

# 3) Populate the freshly inserted row 4 with the new job listing.
$ws.Range("A4").Value = $newTimestamp
$ws.Range("B4").Value = "【急募】WordoressサイトスピードUPのための専門家を探しています!"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5462581"
$ws.Range("G4").Value = 33
$ws.Range("H4").Value = "◇サイト"

# 4) The two pushed-down rows keep their original content, but their
#    fetch timestamp also needs to move to the new run time.
$ws.Range("A5").Value = $newTimestamp
$ws.Range("A6").Value = $newTimestamp

# 5) Row insertion does not re-target the worksheet's <hyperlinks> entries
#    (their `ref` stays pinned to the old row), so rebuild all of them
#    from scratch against the final layout.
$ws.Range("F2").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5217096")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5457458")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5462581")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5462397")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5462522")

# Hyperlinks.Add stamps a fresh cell style for the anchor; put every link
# cell back on the shared "Hyperlink" style so the style table doesn't
# accumulate duplicates.
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("F3").Style = "Hyperlink"
$ws.Range("F4").Style = "Hyperlink"
$ws.Range("F5").Style = "Hyperlink"
$ws.Range("F6").Style = "Hyperlink"

Write-Output "ok"
